# #5: property boat&car done
#
# The 汽車 (car) sheet's row 1 was actually a second copy of the data row
# instead of real column headers, and the sheet was missing the trailing
# property_category / category / date / legislator_name / legislator_id /
# source_file / index columns that every other property sheet has, plus a
# "capacity" (engine displacement) column. Bring it in line with the other
# sheets (using 土地/land's sheet as the template for style + shared
# strings) and append the missing columns to both the header and the data
# row.

$wb = $excel.ActiveWorkbook
$ws3 = $wb.Worksheets.Item(3)   # 汽車 (car)
$ws1 = $wb.Worksheets.Item(1)   # 土地 (land) - used as a formatting/label template

# --- Header row (row 1) --------------------------------------------------
$ws1.Range("B1").Copy($ws3.Range("B1"))          # name
$ws3.Range("C1").Value = "capacity"               # capacity (new column, unique to car sheet)
$ws1.Range("E1:H1").Copy($ws3.Range("D1:G1"))     # owner, register_date, register_reason, acquire_value
$ws1.Range("I1:O1").Copy($ws3.Range("H1:N1"))     # property_category..index

# --- Data row (row 2) -----------------------------------------------------
# A2..G2 (index, name, capacity, owner, register_date, register_reason,
# acquire_value) already hold the right values - only the previously
# missing trailing columns need to be filled in.
$ws1.Range("I2:N2").Copy($ws3.Range("H2:M2"))     # land, normal, 2012-04-30, 鄭天財, 1763, tmp1c9c1
$ws3.Range("N2").Value = 46                        # index (same as column A's row id)
